# Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Albahaca
# Insert 2 new weekly rows of data (rows 529-530), shifting all the
# existing rows 529:577 down to 531:579.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 529, pushing the
# previous rows 529..577 down to 531..579.
$ws.Rows("529:530").Insert()

# Populate the two newly-inserted rows with the new weekly records.
$ws.Range("A529:R530").Value = @(
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44918, 13, 100112052, "Albahaca", "Sin especificar", "Primera", 870, 3000, 4500, 3833, "$/docena de matas", "Región Metropolitana", 639, 6, "Hortaliza"),
    @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44918, 13, 100112052, "Albahaca", "Sin especificar", "Segunda", 350, 2500, 3500, 3071, "$/docena de matas", "Región Metropolitana", 512, 6, "Hortaliza")
)

$ws.Range("D529:D530").NumberFormat = "YYYY-MM-DD HH:MM:SS"
